$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

$ws.Range("D2").Value = '41.178.92'
$ws.Range("E2").Value = '  -1.23%  '
$ws.Range("D3").Value = '2.430.00'
$ws.Range("E3").Value = '  -1.96%  '
$ws.Range("E4").Value = '  +0.09%  '
Set-TextValue $ws.Range("D5") '316.89'
$ws.Range("E5").Value = '  -0.66%  '
Set-TextValue $ws.Range("D6") '88.81'
$ws.Range("E6").Value = '  -4.63%  '
Set-TextValue $ws.Range("D7") '0.542'
$ws.Range("E7").Value = '  -2.41%  '
$ws.Range("E8").Value = '  +0.02%  '
Set-TextValue $ws.Range("D9") '0.496'
$ws.Range("E9").Value = '  -4.51%  '
$ws.Range("B10").Value = 'Avalanche'
$ws.Range("C10").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue $ws.Range("D10") '32.11'
$ws.Range("E10").Value = '  -3.00%  '
$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue $ws.Range("D11") '0.0834'
$ws.Range("E11").Value = '  -5.67%  '
$ws.Range("E12").Value = '  -2.86%  '
$ws.Range("D13").Value = '2.804.58'
$ws.Range("E13").Value = '  -1.94%  '
Set-TextValue $ws.Range("D14") '6.72'
$ws.Range("E14").Value = '  -3.42%  '
Set-TextValue $ws.Range("D15") '15.59'
$ws.Range("E15").Value = '  -0.78%  '
$ws.Range("D16").Value = '2.439.53'
$ws.Range("E16").Value = '  -2.03%  '
Set-TextValue $ws.Range("D17") '0.773'
$ws.Range("E17").Value = '  -2.66%  '
$ws.Range("D18").Value = '41.133.21'
$ws.Range("E18").Value = '  -1.24%  '
$ws.Range("D19").Value = '0.0₃0924'
$ws.Range("E19").Value = '  -4.06%  '
Set-TextValue $ws.Range("D20") '6.24'
$ws.Range("E20").Value = '  -4.15%  '
Set-TextValue $ws.Range("D21") '71.81'
$ws.Range("E21").Value = '  +0.30%  '
Set-TextValue $ws.Range("D22") '11.03'
$ws.Range("E22").Value = '  -4.64%  '
Set-TextValue $ws.Range("D23") '235.74'
$ws.Range("E23").Value = '  -2.68%  '
Set-TextValue $ws.Range("D24") '2.70'
$ws.Range("E24").Value = '  -2.31%  '
$ws.Range("E25").Value = '  +0.06%  '
Set-TextValue $ws.Range("D26") '1.88'
$ws.Range("E26").Value = '  -2.74%  '
Set-TextValue $ws.Range("D27") '23.99'
$ws.Range("E27").Value = '  -3.79%  '
$ws.Range("E28").Value = '  -3.13%  '
Set-TextValue $ws.Range("D29") '9.56'
$ws.Range("E29").Value = '  -3.59%  '
Set-TextValue $ws.Range("D30") '34.63'
$ws.Range("E30").Value = '  -5.63%  '
Set-TextValue $ws.Range("D31") '157.17'
$ws.Range("E31").Value = '  +0.15%  '
$ws.Range("E32").Value = '  -0.02%  '
Set-TextValue $ws.Range("D33") '5.26'
$ws.Range("E33").Value = '  -5.19%  '
Set-TextValue $ws.Range("D34") '2.52'
$ws.Range("E34").Value = '  -2.19%  '
Set-TextValue $ws.Range("D35") '0.0745'
$ws.Range("E35").Value = '  -3.64%  '
Set-TextValue $ws.Range("D36") '2.91'
$ws.Range("E36").Value = '  -0.64%  '
Set-TextValue $ws.Range("D37") '16.57'
$ws.Range("E37").Value = '  -5.74%  '
Set-TextValue $ws.Range("D38") '0.115'
$ws.Range("E38").Value = '  -0.72%  '
Set-TextValue $ws.Range("D39") '1.77'
$ws.Range("E39").Value = '  -3.84%  '
$ws.Range("E40").Value = '  -4.47%  '
Set-TextValue $ws.Range("D41") '3.87'
$ws.Range("E41").Value = '  -4.02%  '
$ws.Range("E42").Value = '  -7.03%  '
$ws.Range("D43").Value = '1.982.34'
$ws.Range("E43").Value = '  -0.15%  '
Set-TextValue $ws.Range("D44") '0.0275'
$ws.Range("E44").Value = '  -4.16%  '
Set-TextValue $ws.Range("D45") '18.30'
$ws.Range("E45").Value = '  -6.96%  '
Set-TextValue $ws.Range("D46") '2.87'
$ws.Range("E46").Value = '  -5.78%  '
Set-TextValue $ws.Range("D47") '9.49'
$ws.Range("E47").Value = '  +2.90%  '
$ws.Range("D48").Value = '2.671.01'
$ws.Range("E48").Value = '  -1.54%  '
Set-TextValue $ws.Range("D49") '95.36'
$ws.Range("E49").Value = '  -2.45%  '
Set-TextValue $ws.Range("D50") '73.30'
$ws.Range("E50").Value = '  -1.24%  '
Set-TextValue $ws.Range("D51") '51.87'
$ws.Range("E51").Value = '  -1.87%  '
